# LA Lakers roster housekeeping edit:
#  1. Swap the Lonnie Walker IV / Austin Reaves rows (rows 5 & 6)
#  2. Move D'Angelo Russell up above Scotty Pippen Jr. (TW) / Cole Swider (TW)
#     (row 16 -> row 14, shifting the other two down one row)
#  3. Correct Malik Beasley's listed height (row 11)
#  4. Correct Davon Reed's listed height and college (row 17)
#
# Row data (columns B:K - the "No." through "bbref url" fields) is moved
# with Range.Copy so that cell data types (numbers vs. shared-string text)
# are preserved exactly, the same way a real Excel row move/sort would
# behave. Column A (the 0-based helper index) is left untouched since it
# never changes in the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch rows far away from the real data, used as temporary holding
# space while rows are rotated/swapped.
$scratch1 = "B100:K100"
$scratch2 = "B101:K101"
$scratch3 = "B102:K102"

# --- 1. Swap row 5 (Lonnie Walker IV) and row 6 (Austin Reaves) ---------
$ws.Range("B5:K5").Copy($ws.Range($scratch1))
$ws.Range("B6:K6").Copy($ws.Range("B5:K5"))
$ws.Range($scratch1).Copy($ws.Range("B6:K6"))
$ws.Range($scratch1).Clear()

# --- 2. Move D'Angelo Russell (row 16) above rows 14-15 -----------------
$ws.Range("B14:K14").Copy($ws.Range($scratch1))
$ws.Range("B15:K15").Copy($ws.Range($scratch2))
$ws.Range("B16:K16").Copy($ws.Range($scratch3))

$ws.Range($scratch3).Copy($ws.Range("B14:K14"))
$ws.Range($scratch1).Copy($ws.Range("B15:K15"))
$ws.Range($scratch2).Copy($ws.Range("B16:K16"))

$ws.Range($scratch1).Clear()
$ws.Range($scratch2).Clear()
$ws.Range($scratch3).Clear()

# --- 3. Correct Malik Beasley's height (row 11) --------------------------
$ws.Range("E11").Value2 = "6-5"

# --- 4. Correct Davon Reed's height and college (row 17) -----------------
$ws.Range("E17").Value2 = "6-4"
$ws.Range("J17").Value2 = "Wichita State, Oklahoma"
